$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.872.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +7.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.312.47"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +6.21%  "

$ws.Range("E4").Value = "  -0.67%  "

$ws.Range("E5").Value = "  +2.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +14.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.573"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.08%  "

$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +12.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.73"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +12.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0799"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.34"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +11.23%  "

$ws.Range("E13").Value = "  +2.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.663.24"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.308.02"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.97"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +10.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.816"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +9.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.760.33"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +7.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.14"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +25.37%  "

$ws.Range("E20").Value = "  +9.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.13"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +7.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.84"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +7.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +10.58%  "

$ws.Range("E24").Value = "  +7.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +11.29%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.86"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +25.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.84"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +9.43%  "

$ws.Range("E30").Value = "  +7.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.75"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +12.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.43"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0799"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +13.15%  "

$ws.Range("E34").Value = "  +6.41%  "

$ws.Range("E35").Value = "  +10.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +13.98%  "

$ws.Range("E37").Value = "  +4.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.79"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +12.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.78"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +23.51%  "

$ws.Range("E40").Value = "  +16.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.43"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +15.64%  "

$ws.Range("E42").Value = "  +12.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +24.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.839.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.81"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +27.15%  "

$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "76.07"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +17.33%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.198"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +19.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.96"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +13.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.21"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.13"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +13.86%  "
